$d = $word.ActiveDocument

# 1. "This contrast with ... memory is accessed is more linear..." ->
#    "This contrasts with ... memory access is more linear..."
$d.Content.Find.Execute(
    " This contrast with regular programs in which memory is accessed is more linear in nature and less computing time is spent navigating the data structure. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " This contrasts with regular programs in which memory access is more linear in nature and less computing time is spent navigating the data structure. ",
    2) | Out-Null

# 2. "that this registers ar" + "e fix" -> "that these registers ar" + "e fix"
$d.Content.Find.Execute(
    " that this registers are fixed",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " that these registers are fixed",
    2) | Out-Null

# 3. remove comma: "...this architecture, in comparison to RISC" -> "...this architecture in comparison to RISC"
$d.Content.Find.Execute(
    " when the C code is compiled for this architecture, in comparison to RISC",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " when the C code is compiled for this architecture in comparison to RISC",
    2) | Out-Null

# 4. "similar RISCV" -> "like RISCV"
$d.Content.Find.Execute(
    "similar RISCV in terms of instruction counts, which is good since fewer instructions count improve power consumption.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "like RISCV in terms of instruction counts, which is good since fewer instructions count improve power consumption.",
    2) | Out-Null

# Merge the following paragraph into this one: drop the paragraph break and
# the leading sentence "The solution presented ClockHands, simulation was
# fundamental to estimate power consumption.", keeping the rest as part of
# the previous paragraph.
$d.Content.Find.Execute(
    "power consumption.^p^tThe solution presented ClockHands, simulation was fundamental to estimate power consumption. Overall,",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "power consumption. Overall,",
    2) | Out-Null

# 5. "We used a cycle-accurate simulator, Onikiri2, for the performance" ->
#    "a cycle-accurate simulator, Onikiri2, for the performance"
$d.Content.Find.Execute(
    "We used a cycle-accurate simulator, Onikiri2, for the performance",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "a cycle-accurate simulator, Onikiri2, for the performance",
    2) | Out-Null
